$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'320.84"
$ws.Range("E2").Value = "'5.58%"
$ws.Range("G2").Value = "'4"
$ws.Range("D3").Value = "'36.17"
$ws.Range("E3").Value = "'0.00%"
$ws.Range("G3").Value = "'4"
$ws.Range("D4").Value = "'5.156"
$ws.Range("E4").Value = "'2.49%"
$ws.Range("G4").Value = "'4"
$ws.Range("D5").Value = "'0.08099"
$ws.Range("E5").Value = "'3.70%"
$ws.Range("G5").Value = "'4"
$ws.Range("D6").Value = "'2.159"
$ws.Range("E6").Value = "'-1.61%"
$ws.Range("G6").Value = "'4"
$ws.Range("D7").Value = "'8.068"
$ws.Range("E7").Value = "'2.07%"
$ws.Range("G7").Value = "'4"
$ws.Range("D8").Value = "'4.129"
$ws.Range("G8").Value = "'4"
$ws.Range("D9").Value = "'0.9285"
$ws.Range("E9").Value = "'1.36%"
$ws.Range("G9").Value = "'4"
$ws.Range("D10").Value = "'0.1010"
$ws.Range("E10").Value = "'4.09%"
$ws.Range("G10").Value = "'4"
$ws.Range("D11").Value = "'0.1883"
$ws.Range("E11").Value = "'1.15%"
$ws.Range("G11").Value = "'4"
$ws.Range("D12").Value = "'0.09204"
$ws.Range("E12").Value = "'7.06%"
$ws.Range("G12").Value = "'4"
$ws.Range("D13").Value = "'0.03571"
$ws.Range("E13").Value = "'2.15%"
$ws.Range("G13").Value = "'4"
$ws.Range("D14").Value = "'0.09934"
$ws.Range("E14").Value = "'0.21%"
$ws.Range("G14").Value = "'4"
$ws.Range("D15").Value = "'0.001436"
$ws.Range("E15").Value = "'0.40%"
$ws.Range("G15").Value = "'4"
$ws.Range("D16").Value = "'0.005650"
$ws.Range("E16").Value = "'-0.41%"
$ws.Range("G16").Value = "'4"
$ws.Range("D17").Value = "'3.455"
$ws.Range("E17").Value = "'-0.09%"
$ws.Range("G17").Value = "'4"
$ws.Range("E18").Value = "'17.96%"
$ws.Range("G18").Value = "'4"
$ws.Range("D19").Value = "'0.3373"
$ws.Range("E19").Value = "'-1.26%"
$ws.Range("G19").Value = "'4"
$ws.Range("D20").Value = "'0.1329"
$ws.Range("E20").Value = "'-1.38%"
$ws.Range("G20").Value = "'4"
$ws.Range("D21").Value = "'5.088"
$ws.Range("E21").Value = "'6.07%"
$ws.Range("G21").Value = "'4"
$ws.Range("D22").Value = "'0.2203"
$ws.Range("E22").Value = "'-0.29%"
$ws.Range("G22").Value = "'4"
$ws.Range("D23").Value = "'0.04606"
$ws.Range("E23").Value = "'-0.30%"
$ws.Range("G23").Value = "'4"
$ws.Range("E24").Value = "'0.91%"
$ws.Range("G24").Value = "'4"
$ws.Range("E25").Value = "'-6.98%"
$ws.Range("G25").Value = "'4"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'-7.19%"
$ws.Range("G26").Value = "'4"
$ws.Range("D27").Value = "'0.0004497"
$ws.Range("E27").Value = "'-5.39%"
$ws.Range("G27").Value = "'4"
$ws.Range("G28").Value = "'4"
$ws.Range("G29").Value = "'4"
$ws.Range("G30").Value = "'4"
$ws.Range("G31").Value = "'4"
$ws.Range("G32").Value = "'4"
$ws.Range("G33").Value = "'4"
$ws.Range("G34").Value = "'4"
$ws.Range("G35").Value = "'4"
$ws.Range("G36").Value = "'4"
$ws.Range("G37").Value = "'4"
$ws.Range("G38").Value = "'4"
$ws.Range("D39").Value = "'0.02013"
$ws.Range("E39").Value = "'10.36%"
$ws.Range("G39").Value = "'4"
$ws.Range("D40").Value = "'0.04984"
$ws.Range("E40").Value = "'5.51%"
$ws.Range("G40").Value = "'4"
$ws.Range("D41").Value = "'0.007816"
$ws.Range("E41").Value = "'4.32%"
$ws.Range("G41").Value = "'4"
$ws.Range("D42").Value = "'0.1402"
$ws.Range("E42").Value = "'0.42%"
$ws.Range("G42").Value = "'4"
$ws.Range("D43").Value = "'0.007823"
$ws.Range("E43").Value = "'0.83%"
$ws.Range("G43").Value = "'4"
$ws.Range("D44").Value = "'0.002081"
$ws.Range("E44").Value = "'-6.73%"
$ws.Range("G44").Value = "'4"
$ws.Range("D45").Value = "'0.01209"
$ws.Range("E45").Value = "'9.04%"
$ws.Range("G45").Value = "'4"
$ws.Range("D46").Value = "'0.00006410"
$ws.Range("E46").Value = "'0.80%"
$ws.Range("G46").Value = "'4"
$ws.Range("E47").Value = "'-0.06%"
$ws.Range("G47").Value = "'4"
$ws.Range("E48").Value = "'14.95%"
$ws.Range("G48").Value = "'4"
$ws.Range("E49").Value = "'-5.13%"
$ws.Range("G49").Value = "'4"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("G50").Value = "'4"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.06%"
$ws.Range("G51").Value = "'4"
